# Refresh the crypto price/volume table on Sheet1 (Price, Volume(1h) columns)
# and correct a handful of rows whose ranking order shifted (Coin / Link swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell -> new literal text, exactly as it should appear in the sheet.
$updates = @(
    @{Cell = "D2"; Value = "62.889.46"}
    @{Cell = "E2"; Value = "  -6.78%  "}
    @{Cell = "D3"; Value = "3.429.23"}
    @{Cell = "E3"; Value = "  -4.31%  "}
    @{Cell = "E4"; Value = "  +0.29%  "}
    @{Cell = "D5"; Value = "384.95"}
    @{Cell = "E5"; Value = "  -7.37%  "}
    @{Cell = "D6"; Value = "121.94"}
    @{Cell = "E6"; Value = "  -5.62%  "}
    @{Cell = "D7"; Value = "3.424.06"}
    @{Cell = "E7"; Value = "  -4.20%  "}
    @{Cell = "D8"; Value = "0.578"}
    @{Cell = "E8"; Value = "  -11.37%  "}
    @{Cell = "E9"; Value = "  +0.12%  "}
    @{Cell = "D10"; Value = "0.656"}
    @{Cell = "E10"; Value = "  -14.92%  "}
    @{Cell = "D11"; Value = "0.139"}
    @{Cell = "E11"; Value = "  -23.66%  "}
    @{Cell = "E12"; Value = "  -14.46%  "}
    @{Cell = "D13"; Value = "38.09"}
    @{Cell = "E13"; Value = "  -10.03%  "}
    @{Cell = "D14"; Value = "3.958.67"}
    @{Cell = "E14"; Value = "  -4.80%  "}
    @{Cell = "D15"; Value = "9.05"}
    @{Cell = "E15"; Value = "  -8.35%  "}
    @{Cell = "E16"; Value = "  -3.18%  "}
    @{Cell = "D17"; Value = "3.420.40"}
    @{Cell = "E17"; Value = "  -4.79%  "}
    @{Cell = "D18"; Value = "18.27"}
    @{Cell = "E18"; Value = "  -10.40%  "}
    @{Cell = "D19"; Value = "12.16"}
    @{Cell = "E19"; Value = "  -0.75%  "}
    @{Cell = "D20"; Value = "62.841.73"}
    @{Cell = "E20"; Value = "  -6.56%  "}
    @{Cell = "D21"; Value = "1.00"}
    @{Cell = "E21"; Value = "  -12.13%  "}
    @{Cell = "D22"; Value = "383.55"}
    @{Cell = "E22"; Value = "  -15.13%  "}
    @{Cell = "D23"; Value = "13.23"}
    @{Cell = "E23"; Value = "  +0.68%  "}
    @{Cell = "D24"; Value = "79.31"}
    @{Cell = "E24"; Value = "  -11.02%  "}
    @{Cell = "D25"; Value = "2.76"}
    @{Cell = "E25"; Value = "  -12.42%  "}
    @{Cell = "D26"; Value = "5.18"}
    @{Cell = "E26"; Value = "  +6.17%  "}
    @{Cell = "D27"; Value = "32.49"}
    @{Cell = "E27"; Value = "  -7.25%  "}
    @{Cell = "D28"; Value = "2.92"}
    @{Cell = "E28"; Value = "  -12.60%  "}
    @{Cell = "D29"; Value = "8.61"}
    @{Cell = "E29"; Value = "  -13.81%  "}
    @{Cell = "D30"; Value = "11.63"}
    @{Cell = "E30"; Value = "  -5.86%  "}
    @{Cell = "D31"; Value = "2.59"}
    @{Cell = "E31"; Value = "  -4.94%  "}
    @{Cell = "E32"; Value = "  -8.48%  "}
    @{Cell = "D33"; Value = "6.50"}
    @{Cell = "E33"; Value = "  -11.94%  "}
    @{Cell = "E34"; Value = "  -8.88%  "}
    @{Cell = "D35"; Value = "0.999"}
    @{Cell = "E35"; Value = "  +0.06%  "}
    @{Cell = "D36"; Value = "36.20"}
    @{Cell = "E36"; Value = "  -11.59%  "}
    @{Cell = "D37"; Value = "53.69"}
    @{Cell = "E37"; Value = "  -5.47%  "}
    @{Cell = "D38"; Value = "0.0425"}
    @{Cell = "E38"; Value = "  -13.99%  "}
    @{Cell = "D39"; Value = "0.992"}
    @{Cell = "E39"; Value = "  -0.56%  "}
    @{Cell = "D40"; Value = "26.39"}
    @{Cell = "E40"; Value = "  +23.70%  "}
    @{Cell = "D41"; Value = "2.61"}
    @{Cell = "E41"; Value = "  +13.13%  "}
    @{Cell = "E42"; Value = "  -12.19%  "}
    @{Cell = "B43"; Value = "ApeXProtocol"}
    @{Cell = "C43"; Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"}
    @{Cell = "D43"; Value = "3.04"}
    @{Cell = "E43"; Value = "  +16.91%  "}
    @{Cell = "B44"; Value = "Monero"}
    @{Cell = "C44"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"}
    @{Cell = "D44"; Value = "136.47"}
    @{Cell = "E44"; Value = "  -8.60%  "}
    @{Cell = "D45"; Value = "0.0₃0586"}
    @{Cell = "E45"; Value = "  -24.74%  "}
    @{Cell = "B46"; Value = "WEMIXToken"}
    @{Cell = "C46"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"}
    @{Cell = "D46"; Value = "2.42"}
    @{Cell = "E46"; Value = "  -11.79%  "}
    @{Cell = "B47"; Value = "ARBITRUM"}
    @{Cell = "C47"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"}
    @{Cell = "D47"; Value = "1.90"}
    @{Cell = "E47"; Value = "  -3.87%  "}
    @{Cell = "B48"; Value = "LidoDAOToken"}
    @{Cell = "C48"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"}
    @{Cell = "D48"; Value = "2.98"}
    @{Cell = "E48"; Value = "  -8.46%  "}
    @{Cell = "B49"; Value = "NEARProtocol"}
    @{Cell = "C49"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"}
    @{Cell = "D49"; Value = "3.92"}
    @{Cell = "E49"; Value = "  -8.96%  "}
    @{Cell = "B50"; Value = "Stacks"}
    @{Cell = "C50"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"}
    @{Cell = "D50"; Value = "2.59"}
    @{Cell = "E50"; Value = "  -15.38%  "}
    @{Cell = "D51"; Value = "0.269"}
    @{Cell = "E51"; Value = "  -14.46%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Value
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a bare number (e.g. "384.95", "1.00"). A plain assignment
        # would make Excel coerce it to a numeric value (losing the significant
        # trailing zero / the text formatting used throughout this column), so
        # force text with a leading quote prefix, then strip the quote-prefix
        # styling it leaves behind so the cell keeps its original (default) style.
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

